$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells whose new value looks like a plain decimal number must be
# forced to Text format first, so Excel stores them as strings rather than
# silently converting to a numeric cell (the sheet keeps all Price values,
# including grouped-thousands ones like "70.295.78", as text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "70.295.78"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "3.596.90"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "581.41"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").Value = "190.07"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "0.631"
$ws.Range("E7").Value = "  -2.23%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.592.98"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("E10").Value = "  +4.19%  "
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").Value = "56.12"
$ws.Range("E12").Value = "  -3.55%  "
$ws.Range("D13").Value = "0.0000311"
$ws.Range("E13").Value = "  +7.55%  "
$ws.Range("D14").Value = "9.70"
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("D15").Value = "4.179.77"
$ws.Range("D16").Value = "19.80"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "3.599.49"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "70.284.09"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").Value = "12.63"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("D22").Value = "490.74"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("D23").Value = "19.55"
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("D24").Value = "4.88"
$ws.Range("E24").Value = "  -9.28%  "
$ws.Range("D25").Value = "96.37"
$ws.Range("E25").Value = "  +6.02%  "
$ws.Range("E26").Value = "  -1.87%  "
$ws.Range("D27").Value = "2.99"
$ws.Range("E27").Value = "  -4.74%  "
$ws.Range("D28").Value = "11.03"
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("D29").Value = "9.40"
$ws.Range("E29").Value = "  -2.74%  "
$ws.Range("D30").Value = "32.25"
$ws.Range("E30").Value = "  -2.33%  "
$ws.Range("E31").Value = "  -3.28%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").Value = "66.15"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.118"
$ws.Range("E34").Value = "  -1.86%  "
$ws.Range("D35").Value = "579.84"
$ws.Range("E35").Value = "  -7.42%  "
$ws.Range("D36").Value = "38.81"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").Value = "0.0₃0814"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").Value = "0.399"
$ws.Range("E39").Value = "  -3.16%  "
$ws.Range("D40").Value = "2.95"
$ws.Range("E40").Value = "  +5.71%  "
$ws.Range("D41").Value = "3.24"
$ws.Range("E41").Value = "  +17.92%  "
$ws.Range("D42").Value = "3.47"
$ws.Range("E42").Value = "  -4.16%  "
$ws.Range("E43").Value = "  -6.40%  "
$ws.Range("D44").Value = "3.222.06"
$ws.Range("E44").Value = "  -2.49%  "
$ws.Range("E45").Value = "  -1.91%  "
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").Value = "9.66"
$ws.Range("E47").Value = "  +6.00%  "
$ws.Range("D48").Value = "3.38"
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "3.22"
$ws.Range("E51").Value = "  -2.87%  "
